$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: D4/E4 were formulas (=I34 / =K34); replace with literal corrected values
$ws.Range("D4").Value = 0.4
$ws.Range("E4").Value = 7.8

# Row 5: D5/E5 were formulas (=I40 / =K40); replace with literal corrected values
$ws.Range("D5").Value = 0.28000000000000003
$ws.Range("E5").Value = 4.9000000000000004

# Row 6: D6/E6 were formulas (=I46 / =K46); replace with literal corrected values
$ws.Range("D6").Value = 0.26
$ws.Range("E6").Value = 4

# Row 7: D7/E7 pointed at the wrong (out-of-range) row 52; fix to point at row 43
$ws.Range("D7").Formula = "=I43"
$ws.Range("E7").Formula = "=K43"

# Move the active selection to I12 (was E12)
[void]$ws.Range("I12").Select()
